$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template rows for copying cell formatting (styles):
#  row 290 -> used when G column has a value (style 1)
#  row 259 -> used when G column is empty (style 2)

# Row 291
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A291:I291").PasteSpecial(-4122) | Out-Null
$ws.Range("A291").Value = 45908
$ws.Range("B291").Value = "Yoan Zouma"
$ws.Range("C291").Value = 55
$ws.Range("D291").Value = 6
$ws.Range("E291").Value = 9
$ws.Range("F291").Value = 4
$ws.Range("G291").Value = "Adducteur cheville "
$ws.Range("H291").Value = 4

# Row 292
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A292:I292").PasteSpecial(-4122) | Out-Null
$ws.Range("A292").Value = 45908
$ws.Range("B292").Value = "Ilyes Boughanmi"
$ws.Range("C292").Value = 55
$ws.Range("D292").Value = 6
$ws.Range("E292").Value = 6
$ws.Range("F292").Value = 4
$ws.Range("G292").Value = "Adducteur "
$ws.Range("H292").Value = 6

# Row 293
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A293:I293").PasteSpecial(-4122) | Out-Null
$ws.Range("A293").Value = 45908
$ws.Range("B293").Value = "Omar Benyounes"
$ws.Range("C293").Value = 55
$ws.Range("D293").Value = 7
$ws.Range("E293").Value = 6
$ws.Range("F293").Value = 0
$ws.Range("H293").Value = 7

# Row 294
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A294:I294").PasteSpecial(-4122) | Out-Null
$ws.Range("A294").Value = 45908
$ws.Range("B294").Value = "Naim Ighbane"
$ws.Range("C294").Value = 55
$ws.Range("D294").Value = 9
$ws.Range("E294").Value = 8
$ws.Range("F294").Value = 0
$ws.Range("H294").Value = 8

# Row 295
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A295:I295").PasteSpecial(-4122) | Out-Null
$ws.Range("A295").Value = 45908
$ws.Range("B295").Value = "Yanis Berrached"
$ws.Range("C295").Value = 55
$ws.Range("D295").Value = 8
$ws.Range("E295").Value = 6
$ws.Range("F295").Value = 0
$ws.Range("H295").Value = 9

# Row 296
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A296:I296").PasteSpecial(-4122) | Out-Null
$ws.Range("A296").Value = 45908
$ws.Range("B296").Value = "Karim Belmahi"
$ws.Range("C296").Value = 55
$ws.Range("D296").Value = 6
$ws.Range("E296").Value = 5
$ws.Range("F296").Value = 0
$ws.Range("H296").Value = 10

# Row 297
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A297:I297").PasteSpecial(-4122) | Out-Null
$ws.Range("A297").Value = 45908
$ws.Range("B297").Value = "Jeremie Laurent"
$ws.Range("C297").Value = 55
$ws.Range("D297").Value = 5
$ws.Range("E297").Value = 5
$ws.Range("F297").Value = 0
$ws.Range("H297").Value = 6

# Row 298
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A298:I298").PasteSpecial(-4122) | Out-Null
$ws.Range("A298").Value = 45908
$ws.Range("B298").Value = "Ilan Ihaddadene"
$ws.Range("C298").Value = 55
$ws.Range("D298").Value = 8
$ws.Range("E298").Value = 7
$ws.Range("F298").Value = 1
$ws.Range("G298").Value = "Semelle "
$ws.Range("H298").Value = 9

# Row 299
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A299:I299").PasteSpecial(-4122) | Out-Null
$ws.Range("A299").Value = 45908
$ws.Range("B299").Value = "Amine Taiar"
$ws.Range("C299").Value = 55
$ws.Range("D299").Value = 7
$ws.Range("E299").Value = 3
$ws.Range("F299").Value = 0
$ws.Range("H299").Value = 3

# Row 300
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A300:I300").PasteSpecial(-4122) | Out-Null
$ws.Range("A300").Value = 45908
$ws.Range("B300").Value = "Emmanuel Valey"
$ws.Range("C300").Value = 55
$ws.Range("D300").Value = 6
$ws.Range("E300").Value = 5
$ws.Range("F300").Value = 1
$ws.Range("G300").Value = "Adducteur "
$ws.Range("H300").Value = 5

# Row 301
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A301:I301").PasteSpecial(-4122) | Out-Null
$ws.Range("A301").Value = 45908
$ws.Range("B301").Value = "Karahali Souaré"
$ws.Range("C301").Value = 55
$ws.Range("D301").Value = 5
$ws.Range("E301").Value = 6
$ws.Range("F301").Value = 6
$ws.Range("G301").Value = "Cheville"
$ws.Range("H301").Value = 7

# Row 302
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A302:I302").PasteSpecial(-4122) | Out-Null
$ws.Range("A302").Value = 45908
$ws.Range("B302").Value = "Naim Dhib"
$ws.Range("C302").Value = 55
$ws.Range("D302").Value = 4
$ws.Range("E302").Value = 8
$ws.Range("F302").Value = 3
$ws.Range("G302").Value = "adducteur"
$ws.Range("H302").Value = 5

# Row 303
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A303:I303").PasteSpecial(-4122) | Out-Null
$ws.Range("A303").Value = 45908
$ws.Range("B303").Value = "Sofiane Belle"
$ws.Range("C303").Value = 55
$ws.Range("D303").Value = 6
$ws.Range("E303").Value = 7
$ws.Range("F303").Value = 2
$ws.Range("G303").Value = "Ischio"
$ws.Range("H303").Value = 5

# Row 304
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A304:I304").PasteSpecial(-4122) | Out-Null
$ws.Range("A304").Value = 45909
$ws.Range("B304").Value = "Amir Etien"
$ws.Range("C304").Value = 70
$ws.Range("D304").Value = 6
$ws.Range("E304").Value = 10
$ws.Range("F304").Value = 0
$ws.Range("H304").Value = 0

# Row 305
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A305:I305").PasteSpecial(-4122) | Out-Null
$ws.Range("A305").Value = 45909
$ws.Range("B305").Value = "Ilyes Boughanmi"
$ws.Range("C305").Value = 70
$ws.Range("D305").Value = 6
$ws.Range("E305").Value = 10
$ws.Range("F305").Value = 5
$ws.Range("G305").Value = "Adducteur "
$ws.Range("H305").Value = 0

# Row 306
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A306:I306").PasteSpecial(-4122) | Out-Null
$ws.Range("A306").Value = 45909
$ws.Range("B306").Value = "Omar Benyounes"
$ws.Range("C306").Value = 70
$ws.Range("D306").Value = 5
$ws.Range("E306").Value = 6
$ws.Range("F306").Value = 3
$ws.Range("G306").Value = "Tibia coup"
$ws.Range("H306").Value = 5

# Row 307
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A307:I307").PasteSpecial(-4122) | Out-Null
$ws.Range("A307").Value = 45909
$ws.Range("B307").Value = "Karim Belmahi"
$ws.Range("C307").Value = 70
$ws.Range("D307").Value = 6
$ws.Range("E307").Value = 8
$ws.Range("F307").Value = 2
$ws.Range("G307").Value = "Courbatures"
$ws.Range("H307").Value = 10

# Row 308
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A308:I308").PasteSpecial(-4122) | Out-Null
$ws.Range("A308").Value = 45909
$ws.Range("B308").Value = "Yoan Zouma"
$ws.Range("C308").Value = 70
$ws.Range("D308").Value = 4
$ws.Range("E308").Value = 5
$ws.Range("F308").Value = 1
$ws.Range("G308").Value = "Malade"
$ws.Range("H308").Value = 0

# Row 309
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A309:I309").PasteSpecial(-4122) | Out-Null
$ws.Range("A309").Value = 45909
$ws.Range("B309").Value = "Yoann Martelat"
$ws.Range("C309").Value = 70
$ws.Range("D309").Value = 6
$ws.Range("E309").Value = 6
$ws.Range("F309").Value = 6
$ws.Range("G309").Value = "Genou"
$ws.Range("H309").Value = 8

# Row 310
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A310:I310").PasteSpecial(-4122) | Out-Null
$ws.Range("A310").Value = 45909
$ws.Range("B310").Value = "Jeremie Laurent"
$ws.Range("C310").Value = 70
$ws.Range("D310").Value = 6
$ws.Range("E310").Value = 6
$ws.Range("F310").Value = 0
$ws.Range("H310").Value = 8

# Row 311
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A311:I311").PasteSpecial(-4122) | Out-Null
$ws.Range("A311").Value = 45909
$ws.Range("B311").Value = "Naim Ighbane"
$ws.Range("C311").Value = 70
$ws.Range("D311").Value = 5
$ws.Range("E311").Value = 6
$ws.Range("F311").Value = 0
$ws.Range("H311").Value = 7

# Row 312
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A312:I312").PasteSpecial(-4122) | Out-Null
$ws.Range("A312").Value = 45909
$ws.Range("B312").Value = "Ilan Ihaddadene"
$ws.Range("C312").Value = 70
$ws.Range("D312").Value = 6
$ws.Range("E312").Value = 6
$ws.Range("F312").Value = 0
$ws.Range("H312").Value = 7

# Row 313
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A313:I313").PasteSpecial(-4122) | Out-Null
$ws.Range("A313").Value = 45909
$ws.Range("B313").Value = "Emmanuel Valey"
$ws.Range("C313").Value = 70
$ws.Range("D313").Value = 6
$ws.Range("E313").Value = 7
$ws.Range("F313").Value = 1
$ws.Range("G313").Value = "Adducteur "
$ws.Range("H313").Value = 8

# Row 314
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A314:I314").PasteSpecial(-4122) | Out-Null
$ws.Range("A314").Value = 45909
$ws.Range("B314").Value = "Naim Dhib"
$ws.Range("C314").Value = 70
$ws.Range("D314").Value = 5
$ws.Range("E314").Value = 6
$ws.Range("F314").Value = 3
$ws.Range("G314").Value = "Adducteur "
$ws.Range("H314").Value = 7

# Row 315
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A315:I315").PasteSpecial(-4122) | Out-Null
$ws.Range("A315").Value = 45909
$ws.Range("B315").Value = "Karahali Souaré"
$ws.Range("C315").Value = 70
$ws.Range("D315").Value = 5
$ws.Range("E315").Value = 6
$ws.Range("F315").Value = 6
$ws.Range("G315").Value = "Cheville"
$ws.Range("H315").Value = 6

# Row 316
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A316:I316").PasteSpecial(-4122) | Out-Null
$ws.Range("A316").Value = 45909
$ws.Range("B316").Value = "Mattheo Haon"
$ws.Range("C316").Value = 70
$ws.Range("D316").Value = 6
$ws.Range("E316").Value = 6
$ws.Range("F316").Value = 4
$ws.Range("G316").Value = "Adducteur "
$ws.Range("H316").Value = 8

# Row 317
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A317:I317").PasteSpecial(-4122) | Out-Null
$ws.Range("A317").Value = 45909
$ws.Range("B317").Value = "Wael Fareh"
$ws.Range("C317").Value = 70
$ws.Range("D317").Value = 5
$ws.Range("E317").Value = 2
$ws.Range("F317").Value = 0
$ws.Range("H317").Value = 6

# Row 318
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A318:I318").PasteSpecial(-4122) | Out-Null
$ws.Range("A318").Value = 45909
$ws.Range("B318").Value = "Hedi Nasri"
$ws.Range("C318").Value = 70
$ws.Range("D318").Value = 5
$ws.Range("E318").Value = 2
$ws.Range("F318").Value = 0
$ws.Range("H318").Value = 7

# Row 319
$ws.Range("A290:I290").Copy() | Out-Null
$ws.Range("A319:I319").PasteSpecial(-4122) | Out-Null
$ws.Range("A319").Value = 45909
$ws.Range("B319").Value = "Levy Ndoutoume"
$ws.Range("C319").Value = 70
$ws.Range("D319").Value = 6
$ws.Range("E319").Value = 7
$ws.Range("F319").Value = 3
$ws.Range("G319").Value = "Cheville ischio"
$ws.Range("H319").Value = 7

# Set shared formula for I291:I319 (adjusts relative refs automatically)
$ws.Range("I291:I319").Formula = "=C291*D291"

# Update selection to match target view state
$ws.Range("K310").Select() | Out-Null
